$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data
#    rows (2 through 424) from 45188 (2023-09-19) to 45189 (2023-09-20).
$ws.Range("C2:C424").Value2 = 45189

# 2. Row 424 previously had no explicit row height; now it gets one
#    (it is no longer the last row in the sheet).
$ws.Rows.Item(424).RowHeight = 15

# 3. Append the two new data rows (425 and 426).
$ws.Range("A425").Value2 = "A 44201-2023"
$ws.Range("B425").Value2 = 45188
$ws.Range("C425").Value2 = 45189
$ws.Range("D425").Value2 = "DALARNAS LÄN"
$ws.Range("E425").Value2 = "MORA"
$ws.Range("G425").Value2 = 0.6
$ws.Range("H425").Value2 = 0
$ws.Range("I425").Value2 = 0
$ws.Range("J425").Value2 = 0
$ws.Range("K425").Value2 = 0
$ws.Range("L425").Value2 = 0
$ws.Range("M425").Value2 = 0
$ws.Range("N425").Value2 = 0
$ws.Range("O425").Value2 = 0
$ws.Range("P425").Value2 = 0
$ws.Range("Q425").Value2 = 0

$ws.Range("A426").Value2 = "A 44205-2023"
$ws.Range("B426").Value2 = 45188
$ws.Range("C426").Value2 = 45189
$ws.Range("D426").Value2 = "DALARNAS LÄN"
$ws.Range("E426").Value2 = "MORA"
$ws.Range("G426").Value2 = 0.7
$ws.Range("H426").Value2 = 0
$ws.Range("I426").Value2 = 0
$ws.Range("J426").Value2 = 0
$ws.Range("K426").Value2 = 0
$ws.Range("L426").Value2 = 0
$ws.Range("M426").Value2 = 0
$ws.Range("N426").Value2 = 0
$ws.Range("O426").Value2 = 0
$ws.Range("P426").Value2 = 0
$ws.Range("Q426").Value2 = 0

# Apply the date number format to the new B/C date cells (reuses the
# existing "YYYY-MM-DD" style, same as the rest of the column).
$ws.Range("B425:C426").NumberFormat = "YYYY-MM-DD"

# Apply the wrap-text style to the (empty) R cells of the new rows, as
# done for every other row in the sheet.
$ws.Range("R425").WrapText = $true
$ws.Range("R426").WrapText = $true

# Row 425 is not the last row any more, so it gets an explicit row
# height like every other non-final row. Row 426, being the new last
# row, keeps the implicit (non-custom) height.
$ws.Rows.Item(425).RowHeight = 15
